# Update the ALNS staff scheduling roster (Solution sheet, B2:AC10)
# with the new requirement values produced by the ADR-2614 change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Solution")

# Index 0 corresponds to worksheet row 2 (Staff_1), index 8 to row 10 (Staff_9).
$scheduleData = @(
    @("A1", "A1", "M3", "A1", "DO", "M1", "M1", "M3", "A1", "M1", "DO", "M1", "M1", "M1", "A1", "A1", "M1", "A1", "DO", "M3", "M1", "A1", "A1", "M3", "DO", "M1", "M1", "M1"),
    @("M2", "M2", "A2", "A2", "DO", "M1", "A1", "M2", "DO", "M2", "A2", "A2", "M1", "A1", "M2", "PH", "PH", "DO", "M2", "M2", "A2", "M2", "PH", "PH", "M2", "A2", "DO", "M2"),
    @("DO", "M1", "M1", "M1", "A1", "A1", "M3", "M1", "M3", "A1", "DO", "M1", "M1", "M1", "M1", "PH", "PH", "M3", "A1", "M1", "DO", "A1", "PH", "PH", "M1", "DO", "M3", "A1"),
    @("DO", "M2", "M2", "M1", "M1", "M2", "M2", "A2", "M2", "DO", "M1", "M2", "A1", "A2", "DO", "PH", "PH", "M2", "M2", "A2", "M2", "M2", "PH", "PH", "DO", "M2", "M2", "M2"),
    @("M2", "M2", "A1", "DO", "M2", "M2", "A1", "M2", "M1", "A2", "M2", "DO", "M2", "A1", "M2", "M2", "A1", "A2", "DO", "M2", "M1", "M2", "M2", "A2", "A1", "DO", "M2", "A1"),
    @("A1", "A1", "M3", "M1", "A1", "A1", "DO", "A1", "A1", "M3", "A1", "A1", "M1", "DO", "A1", "A1", "M1", "DO", "M3", "A1", "A1", "A1", "A1", "M1", "M3", "A1", "A1", "DO"),
    @("M2", "A2", "A1", "M1", "A2", "M2", "DO", "A2", "A2", "M2", "A1", "A1", "A2", "DO", "A2", "A2", "M2", "A1", "A1", "M2", "DO", "DO", "M1", "A2", "A2", "M2", "M1", "A2"),
    @("M2", "M2", "A2", "DO", "M1", "A2", "A1", "A2", "DO", "M2", "M2", "M2", "M1", "M1", "DO", "M2", "A2", "M2", "M1", "M1", "A2", "M2", "A2", "M2", "A2", "M1", "A1", "DO"),
    @("A2", "A2", "M2", "A2", "M1", "M1", "DO", "DO", "M1", "A2", "M1", "M2", "A2", "A2", "M2", "PH", "PH", "M2", "A2", "A2", "DO", "A2", "PH", "PH", "M2", "A2", "A2", "DO")
)

for ($r = 0; $r -lt $scheduleData.Length; $r++) {
    $rowValues = $scheduleData[$r]
    $rowNum = $r + 2
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $ws.Cells.Item($rowNum, $c + 2).Value = $rowValues[$c]
    }
}

